# Power Consumption.xlsx — add "Flight Data" sheet with mission power &
# energy-consumed-in-cruise calculations, as the 5th (last) worksheet.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the current last sheet so it lands at the
# end of the tab strip (Worksheets.Add() with no args inserts *before* the
# active sheet, which is not what we want here).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Flight Data"

# Header row
$ws.Range("A1").Value = "Part"
$ws.Range("B1").Value = "Time (s)"
$ws.Range("C1").Value = "Average Vertical Thrust (kgf)"
$ws.Range("D1").Value = "Average Horizontal Thrust (kgf)"
$ws.Range("E1").Value = "Average Total Thrust (kgf)"
$ws.Range("F1").Value = "Average Power Consumed (W)"
$ws.Range("G1").Value = "Total Energy Consumed (Wh)"

# Mission-leg labels
$ws.Range("A2").Value = "Take-off"
$ws.Range("A3").Value = "Cruise to SZ"
$ws.Range("A4").Value = "Spray Zone"
$ws.Range("A5").Value = "Cruise to Landing"
$ws.Range("A6").Value = "Landing"

# Cruise-to-spray-zone leg measurements
$ws.Range("B3").Value = 134
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 3.23
$ws.Range("F3").Value = 451

# Cruise-to-landing leg measurements
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 3.23
$ws.Range("F5").Value = 451

# Average total thrust = sqrt(horizontal^2 + vertical^2)
$ws.Range("E2").Formula = "=SQRT(D2*D2+C2*C2)"
$ws.Range("E3").Formula = "=SQRT(D3*D3+C3*C3)"
$ws.Range("E4:E6").Formula = "=SQRT(D4*D4+C4*C4)"

# Total energy consumed (Wh) = time (s) * power (W) / 3600
$ws.Range("G2").Formula = "=B2*F2/3600"
$ws.Range("G3:G6").Formula = "=B3*F3/3600"

# Column widths to fit the long headers
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(2).ColumnWidth = 8.5
$ws.Columns.Item(3).ColumnWidth = 24
$ws.Columns.Item(4).ColumnWidth = 27.5
$ws.Columns.Item(5).ColumnWidth = 23.16666666666667
$ws.Columns.Item(6).ColumnWidth = 26.66666666666667
$ws.Columns.Item(7).ColumnWidth = 27.16666666666667

# Match the author's last selection on this (now active) sheet
[void]$ws.Range("G9").Select()
